$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-15 Saturday" "2025-02-16 Sunday"

Replace-Text "78×80=6240" "86×12=1032"
Replace-Text "81×46=3726" "88×88=7744"
Replace-Text "69×43=2967" "83×29=2407"
Replace-Text "32×19=608" "61×36=2196"
Replace-Text "68×81=5508" "76×87=6612"

Replace-Text "91×59=5369" "75×78=5850"
Replace-Text "90×81=7290" "46×25=1150"
Replace-Text "52×11=572" "71×25=1775"
Replace-Text "14×63=882" "86×21=1806"
Replace-Text "91×62=5642" "26×33=858"

Replace-Text "72×41=2952" "21×89=1869"
Replace-Text "80×70=5600" "12×73=876"
Replace-Text "65×21=1365" "57×42=2394"
Replace-Text "53×36=1908" "32×60=1920"
Replace-Text "95×73=6935" "15×11=165"

Replace-Text "36×56=2016" "27×11=297"
Replace-Text "65×12=780" "56×14=784"
Replace-Text "88×38=3344" "26×76=1976"
Replace-Text "42×86=3612" "37×57=2109"
Replace-Text "86×78=6708" "11×27=297"

Replace-Text "25×24=600" "19×65=1235"
Replace-Text "90×66=5940" "27×94=2538"
Replace-Text "44×39=1716" "23×63=1449"
Replace-Text "30×19=570" "17×11=187"
Replace-Text "84×25=2100" "84×48=4032"
